# "Estimativa de Tamanho" update: refresh TEMPO REAL (D) / recompute DIFERENCA (E),
# a few TEMPO ESTIM. (C) corrections, drop the RFS21/RFS22 rows, and move the
# TOTAL row up so it sums C3:C22 / D3:D22 / E3:E22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Página1")

# --- 1) Turn the current row 23 (still "[RFS21]") into the future TOTAL row ---
# Grab the formatting (styles) of the current TOTAL row (row 25) and stamp it
# onto row 23 first, so the bold/bordered TOTAL style indices get reused
# instead of new styles being minted.
$ws.Range("B25:E25").Copy($ws.Range("B23:E23"))
$ws.Range("B23").Value = "TOTAL"

# --- 2) Remove the old row 24 ([RFS22]) and old row 25 (previous TOTAL line) ---
# This shifts the trailing blank rows up by two (so the sheet ends at row 29
# instead of row 31) and leaves our new TOTAL row sitting at row 23.
$ws.Rows("24:25").Delete()

# --- 3) Fill in TEMPO REAL (column D) for rows that now have a real measurement,
#     copying the green "estimate" style from column C so the fill/border match. ---
$dCellsToFormat = @(3, 4, 5, 7, 8, 9, 10)
foreach ($r in $dCellsToFormat) {
    $ws.Range("C$r").Copy($ws.Range("D$r"))
}

# --- 4) Update TEMPO ESTIM. (C) and TEMPO REAL (D) values per row. ---
$ws.Range("C3").Value = 3.0
$ws.Range("D3").Value = 2.0

$ws.Range("C4").Value = 1.5
$ws.Range("D4").Value = 1.0

$ws.Range("D5").Value = 1.0

$ws.Range("C6").Value = 1.0

$ws.Range("C7").Value = 2.0
$ws.Range("D7").Value = 2.0

$ws.Range("D8").Value = 1.0

$ws.Range("C9").Value = 1.0
$ws.Range("D9").Value = 1.0

$ws.Range("C10").Value = 1.0
$ws.Range("D10").Value = 1.0

$ws.Range("C11").Value = 2.0

$ws.Range("C13").Value = 2.0

$ws.Range("C15").Value = 3.0

$ws.Range("C21").Value = 2.0

$ws.Range("C22").Value = 3.0

# --- 5) Re-establish the shared formulas spanning the (now shorter) data range ---
# DIFERENCA = TEMPO REAL - TEMPO ESTIM. for every activity row ...
$ws.Range("E3:E22").Formula = "=D3-C3"
# ... and the TOTAL row sums each column across the remaining activity rows.
$ws.Range("C23:E23").Formula = "=SUM(C3:C22)"
